$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 829.42224
$ws.Cells.Item(15, 9).Value = 829.42224
$ws.Cells.Item(15, 11).Value = 2488.26672
$ws.Cells.Item(15, 13).Value = -2319.26672
$ws.Cells.Item(49, 8).Value = 400
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).ClearContents()
$ws.Cells.Item(74, 8).Value = 10326
$ws.Cells.Item(74, 9).Value = 8586.857
$ws.Cells.Item(74, 11).Value = 8586.857
$ws.Cells.Item(74, 13).Value = -7650.857
$ws.Cells.Item(77, 8).Value = 10326
$ws.Cells.Item(77, 9).Value = 8586.857
$ws.Cells.Item(77, 11).Value = 42934.285
$ws.Cells.Item(77, 13).Value = -38254.285
$ws.Cells.Item(98, 8).Value = 2798.2307
$ws.Cells.Item(98, 9).Value = 2576.138
$ws.Cells.Item(98, 11).Value = 2576.138
$ws.Cells.Item(98, 13).Value = -1078.138
$ws.Cells.Item(99, 8).Value = 477.25
$ws.Cells.Item(99, 10).Value = 709.5
$ws.Cells.Item(99, 12).Value = 2128.5
$ws.Cells.Item(99, 14).Value = -5124.5
$ws.Cells.Item(100, 8).Value = 9691.412
$ws.Cells.Item(100, 9).Value = 1327.3334
$ws.Cells.Item(100, 11).Value = 1327.3334
$ws.Cells.Item(100, 13).Value = -786.3334
$ws.Cells.Item(122, 8).Value = 2798.2307
$ws.Cells.Item(122, 9).Value = 2576.138
$ws.Cells.Item(122, 11).Value = 7728.414
$ws.Cells.Item(122, 13).Value = -5278.414
$ws.Cells.Item(125, 8).Value = 13893670
$ws.Cells.Item(125, 10).Value = 15878164
$ws.Cells.Item(125, 12).Value = 142903476
$ws.Cells.Item(125, 14).Value = -142908396
$ws.Cells.Item(132, 8).Value = 3131.3809
$ws.Cells.Item(132, 9).Value = 3383.1052
$ws.Cells.Item(132, 11).Value = 10149.3156
$ws.Cells.Item(132, 13).Value = -7619.3156
$ws.Cells.Item(138, 8).Value = 4294.2856
$ws.Cells.Item(138, 10).Value = 5428.816
$ws.Cells.Item(138, 12).Value = 16286.448
$ws.Cells.Item(138, 14).Value = -26566.448

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 60360.65
$ws.Cells.Item(2, 9).Value = 63820.688
$ws.Cells.Item(2, 11).Value = 63820.688
$ws.Cells.Item(2, 13).Value = -63707.688
$ws.Cells.Item(32, 8).Value = 3047.845
$ws.Cells.Item(32, 9).Value = 3089.2856
$ws.Cells.Item(32, 11).Value = 3089.2856
$ws.Cells.Item(32, 13).Value = -2802.2856
$ws.Cells.Item(74, 8).Value = 1824
$ws.Cells.Item(74, 9).Value = 1824
$ws.Cells.Item(74, 11).Value = 1824
$ws.Cells.Item(74, 13).Value = -950
$ws.Cells.Item(77, 8).Value = 1824
$ws.Cells.Item(77, 9).Value = 1824
$ws.Cells.Item(77, 11).Value = 9120
$ws.Cells.Item(77, 13).Value = -4752
$ws.Cells.Item(116, 8).Value = 60360.65
$ws.Cells.Item(116, 9).Value = 63820.688
$ws.Cells.Item(116, 11).Value = 63820.688
$ws.Cells.Item(116, 13).Value = -61526.688

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 60360.65
$ws.Cells.Item(3, 9).Value = 63820.688
$ws.Cells.Item(3, 11).Value = 63820.688
$ws.Cells.Item(3, 13).Value = -63706.688
$ws.Cells.Item(14, 8).Value = 1000
$ws.Cells.Item(14, 10).Value = 1000
$ws.Cells.Item(14, 12).Value = 1000
$ws.Cells.Item(14, 14).Value = -1344
$ws.Cells.Item(80, 8).Value = 741.25
$ws.Cells.Item(80, 10).Value = 932.3333
$ws.Cells.Item(80, 12).Value = 932.3333
$ws.Cells.Item(80, 14).Value = -2928.3333
$ws.Cells.Item(83, 8).Value = 741.25
$ws.Cells.Item(83, 10).Value = 932.3333
$ws.Cells.Item(83, 12).Value = 4661.6665
$ws.Cells.Item(83, 14).Value = -14645.6665

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 35232.562
$ws.Cells.Item(31, 9).Value = 2760
$ws.Cells.Item(31, 10).Value = 42726.23
$ws.Cells.Item(31, 11).Value = 2760
$ws.Cells.Item(31, 12).Value = 42726.23
$ws.Cells.Item(31, 13).Value = -2465
$ws.Cells.Item(31, 14).Value = -43316.23
$ws.Cells.Item(34, 8).Value = 35232.562
$ws.Cells.Item(34, 9).Value = 2760
$ws.Cells.Item(34, 10).Value = 42726.23
$ws.Cells.Item(34, 11).Value = 2760
$ws.Cells.Item(34, 12).Value = 42726.23
$ws.Cells.Item(34, 13).Value = -2558
$ws.Cells.Item(34, 14).Value = -43130.23
$ws.Cells.Item(58, 8).Value = 4928.206
$ws.Cells.Item(58, 9).Value = 4553.979
$ws.Cells.Item(58, 10).Value = 5826.35
$ws.Cells.Item(58, 11).Value = 4553.979
$ws.Cells.Item(58, 12).Value = 5826.35
$ws.Cells.Item(58, 13).Value = -4350.979
$ws.Cells.Item(58, 14).Value = -6232.35
$ws.Cells.Item(136, 8).Value = 4928.206
$ws.Cells.Item(136, 9).Value = 4553.979
$ws.Cells.Item(136, 10).Value = 5826.35
$ws.Cells.Item(136, 11).Value = 13661.937
$ws.Cells.Item(136, 12).Value = 17479.05
$ws.Cells.Item(136, 13).Value = -11111.937
$ws.Cells.Item(136, 14).Value = -22579.05

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2398.7
$ws.Cells.Item(5, 9).Value = 1529
$ws.Cells.Item(5, 10).Value = 2867
$ws.Cells.Item(5, 11).Value = 4587
$ws.Cells.Item(5, 12).Value = 8601
$ws.Cells.Item(5, 13).Value = -4475
$ws.Cells.Item(5, 14).Value = -8825
$ws.Cells.Item(6, 8).Value = 157.04762
$ws.Cells.Item(6, 10).Value = 504.5
$ws.Cells.Item(6, 12).Value = 1513.5
$ws.Cells.Item(6, 14).Value = -1739.5
$ws.Cells.Item(103, 8).Value = 3821.5557
$ws.Cells.Item(103, 10).Value = 6736.8
$ws.Cells.Item(103, 12).Value = 20210.4
$ws.Cells.Item(103, 14).Value = -21968.4
$ws.Cells.Item(107, 8).Value = 94109.23
$ws.Cells.Item(107, 9).Value = 1028
$ws.Cells.Item(107, 11).Value = 3084
$ws.Cells.Item(107, 13).Value = -1164
$ws.Cells.Item(124, 8).Value = 1207.75
$ws.Cells.Item(124, 10).Value = 1999
$ws.Cells.Item(124, 12).Value = 5997
$ws.Cells.Item(124, 14).Value = -15817
$ws.Cells.Item(131, 8).Value = 3115.5483
$ws.Cells.Item(131, 9).Value = 1345.1111
$ws.Cells.Item(131, 11).Value = 4035.3333
$ws.Cells.Item(131, 13).Value = 1004.6667
$ws.Cells.Item(132, 8).Value = 4536.467
$ws.Cells.Item(132, 9).Value = 1616.5
$ws.Cells.Item(132, 10).Value = 5598.273
$ws.Cells.Item(132, 11).Value = 14548.5
$ws.Cells.Item(132, 12).Value = 50384.457
$ws.Cells.Item(132, 13).Value = -12018.5
$ws.Cells.Item(132, 14).Value = -55444.457
$ws.Cells.Item(135, 8).Value = 2398.7
$ws.Cells.Item(135, 9).Value = 1529
$ws.Cells.Item(135, 10).Value = 2867
$ws.Cells.Item(135, 11).Value = 13761
$ws.Cells.Item(135, 12).Value = 25803
$ws.Cells.Item(135, 13).Value = -11226
$ws.Cells.Item(135, 14).Value = -30873

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 803
$ws.Cells.Item(102, 9).Value = 582.4706
$ws.Cells.Item(102, 11).Value = 582.4706
$ws.Cells.Item(102, 13).Value = 1039.5294
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 80000
$ws.Cells.Item(140, 10).Value = 80000
$ws.Cells.Item(140, 12).Value = 80000
$ws.Cells.Item(140, 14).Value = -90360

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 631043.25
$ws.Cells.Item(7, 9).Value = 7927.4287
$ws.Cells.Item(7, 11).Value = 7927.4287
$ws.Cells.Item(7, 13).Value = -7815.4287
$ws.Cells.Item(40, 8).Value = 96610
$ws.Cells.Item(40, 10).Value = 7833.3335
$ws.Cells.Item(40, 12).Value = 7833.3335
$ws.Cells.Item(40, 14).Value = -8105.3335
$ws.Cells.Item(126, 8).Value = 631043.25
$ws.Cells.Item(126, 9).Value = 7927.4287
$ws.Cells.Item(126, 11).Value = 23782.2861
$ws.Cells.Item(126, 13).Value = -21312.2861
$ws.Cells.Item(132, 8).Value = 9183.75
$ws.Cells.Item(132, 9).Value = 6745.5
$ws.Cells.Item(132, 11).Value = 20236.5
$ws.Cells.Item(132, 13).Value = -17706.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 89308
$ws.Cells.Item(62, 9).Value = 204339.2
$ws.Cells.Item(62, 11).Value = 204339.2
$ws.Cells.Item(62, 13).Value = -203715.2
$ws.Cells.Item(65, 8).Value = 89308
$ws.Cells.Item(65, 9).Value = 204339.2
$ws.Cells.Item(65, 11).Value = 1021696
$ws.Cells.Item(65, 13).Value = -1018576
$ws.Cells.Item(107, 8).Value = 725.24
$ws.Cells.Item(107, 10).Value = 429.625
$ws.Cells.Item(107, 12).Value = 1288.875
$ws.Cells.Item(107, 14).Value = -5128.875
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).ClearContents()
